$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170..234 down to 171..235.
$ws.Rows(170).Insert()

# Fill in the newly inserted row 170 with the new record's data.
# Columns that duplicate the surrounding records (A,B,C,E,F,G,H,I,J,L,Q,R,T) are
# filled in too, matching the rest of the "Naranja" block.
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44468
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100102
$ws.Range("H170").Value = "Cítricos"
$ws.Range("I170").Value = 100102005
$ws.Range("J170").Value = "Naranja"
$ws.Range("K170").Value = "Lane Late"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 200
$ws.Range("N170").Value = 13000
$ws.Range("O170").Value = 13000
$ws.Range("P170").Value = 13000
$ws.Range("Q170").Value = "$/caja 15 kilos empedrada"
$ws.Range("R170").Value = "Región de O'Higgins"
$ws.Range("S170").Value = 867
$ws.Range("T170").Value = 15
